$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the columns we are about to touch so that
# Excel does not reinterpret numeric-looking / percent-looking strings
# (e.g. "288.11", "-1.01%", "8") as numbers or percentages.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply the updated values scraped on Fri Jan 20 08:14:31 UTC 2023
$ws.Range("D2").Value = "288.11"
$ws.Range("E2").Value = "-1.01%"
$ws.Range("G2").Value = "8"
$ws.Range("D3").Value = "31.04"
$ws.Range("E3").Value = "1.26%"
$ws.Range("G3").Value = "8"
$ws.Range("D4").Value = "4.918"
$ws.Range("E4").Value = "-0.79%"
$ws.Range("G4").Value = "8"
$ws.Range("E5").Value = "1.65%"
$ws.Range("G5").Value = "8"
$ws.Range("D6").Value = "2.246"
$ws.Range("E6").Value = "26.00%"
$ws.Range("G6").Value = "8"
$ws.Range("E7").Value = "0.47%"
$ws.Range("G7").Value = "8"
$ws.Range("D8").Value = "3.733"
$ws.Range("E8").Value = "-0.70%"
$ws.Range("G8").Value = "8"
$ws.Range("D9").Value = "0.9033"
$ws.Range("E9").Value = "0.74%"
$ws.Range("G9").Value = "8"
$ws.Range("D10").Value = "0.09319"
$ws.Range("E10").Value = "20.84%"
$ws.Range("G10").Value = "8"
$ws.Range("E11").Value = "2.04%"
$ws.Range("G11").Value = "8"
$ws.Range("D12").Value = "0.08225"
$ws.Range("E12").Value = "1.81%"
$ws.Range("G12").Value = "8"
$ws.Range("D13").Value = "0.03121"
$ws.Range("E13").Value = "2.84%"
$ws.Range("G13").Value = "8"
$ws.Range("D14").Value = "0.09945"
$ws.Range("E14").Value = "-0.79%"
$ws.Range("G14").Value = "8"
$ws.Range("D15").Value = "0.001495"
$ws.Range("E15").Value = "-0.21%"
$ws.Range("G15").Value = "8"
$ws.Range("D16").Value = "0.005710"
$ws.Range("E16").Value = "-2.68%"
$ws.Range("G16").Value = "8"
$ws.Range("D17").Value = "3.495"
$ws.Range("E17").Value = "0.79%"
$ws.Range("G17").Value = "8"
$ws.Range("D18").Value = "2.059"
$ws.Range("E18").Value = "-1.21%"
$ws.Range("G18").Value = "8"
$ws.Range("D19").Value = "0.3331"
$ws.Range("E19").Value = "0.43%"
$ws.Range("G19").Value = "8"
$ws.Range("E20").Value = "-0.01%"
$ws.Range("G20").Value = "8"
$ws.Range("D21").Value = "4.165"
$ws.Range("E21").Value = "3.27%"
$ws.Range("G21").Value = "8"
$ws.Range("D22").Value = "0.2112"
$ws.Range("E22").Value = "-6.16%"
$ws.Range("G22").Value = "8"
$ws.Range("D23").Value = "0.04521"
$ws.Range("E23").Value = "0.27%"
$ws.Range("G23").Value = "8"
$ws.Range("E24").Value = "-0.40%"
$ws.Range("G24").Value = "8"
$ws.Range("D25").Value = "0.004155"
$ws.Range("E25").Value = "3.51%"
$ws.Range("G25").Value = "8"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "3.99%"
$ws.Range("G26").Value = "8"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("E27").Value = "-95.47%"
$ws.Range("G27").Value = "8"
$ws.Range("G28").Value = "8"
$ws.Range("G29").Value = "8"
$ws.Range("G30").Value = "8"
$ws.Range("G31").Value = "8"
$ws.Range("G32").Value = "8"
$ws.Range("G33").Value = "8"
$ws.Range("G34").Value = "8"
$ws.Range("G35").Value = "8"
$ws.Range("G36").Value = "8"
$ws.Range("G37").Value = "8"
$ws.Range("G38").Value = "8"
$ws.Range("D39").Value = "0.01568"
$ws.Range("E39").Value = "-2.27%"
$ws.Range("G39").Value = "8"
$ws.Range("D40").Value = "0.04446"
$ws.Range("E40").Value = "0.94%"
$ws.Range("G40").Value = "8"
$ws.Range("D41").Value = "0.007338"
$ws.Range("E41").Value = "0.75%"
$ws.Range("G41").Value = "8"
$ws.Range("D42").Value = "0.009568"
$ws.Range("E42").Value = "24.52%"
$ws.Range("G42").Value = "8"
$ws.Range("D43").Value = "0.1325"
$ws.Range("E43").Value = "1.16%"
$ws.Range("G43").Value = "8"
$ws.Range("D44").Value = "0.002250"
$ws.Range("E44").Value = "18.41%"
$ws.Range("G44").Value = "8"
$ws.Range("D45").Value = "0.009101"
$ws.Range("E45").Value = "-4.33%"
$ws.Range("G45").Value = "8"
$ws.Range("D46").Value = "0.00006125"
$ws.Range("E46").Value = "2.92%"
$ws.Range("G46").Value = "8"
$ws.Range("E47").Value = "0.00%"
$ws.Range("G47").Value = "8"
$ws.Range("D48").Value = "2.258"
$ws.Range("E48").Value = "0.54%"
$ws.Range("G48").Value = "8"
$ws.Range("D49").Value = "0.002000"
$ws.Range("E49").Value = "-33.34%"
$ws.Range("G49").Value = "8"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "0.00%"
$ws.Range("G50").Value = "8"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "0.00%"
$ws.Range("G51").Value = "8"
